$wb = $excel.ActiveWorkbook

# Sheet1: survey
$ws1 = $wb.Worksheets.Item("survey")
$ws1.Range("F2").Value = "plot_id"
$ws1.Range("G4").Value = "What is the height of the plant in inches?"

# Sheet3: queries
$ws3 = $wb.Worksheets.Item("queries")
$ws3.Range("E2").Value = "_id >= ?"

# Sheet4: settings
$ws4 = $wb.Worksheets.Item("settings")
$ws4.Range("A5").Value = "table_id"
$ws4.Range("B5").Value = "visit"
$ws4.Range("A6").Value = "disableSwipeNavigation"
$ws4.Range("B6").Value = $true

$ws1.Activate()
